$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Fishing)
$ws.Range("C2").Value = 0.5437913689595455
$ws.Range("D2").Value = 0.5437913689595455
$ws.Range("E2").Value = 2.2314014367748425
$ws.Range("F2").Value = 0.015053500237779254
$ws.Range("G2").Value = 0.0329

# Row 3 (Residuals)
$ws.Range("C3").Value = 35.580124024140254
$ws.Range("D3").Value = 0.243699479617399
$ws.Range("F3").Value = 0.9849464997622207

# Row 4 (Total)
$ws.Range("C4").Value = 36.1239153930998
